# Add the Zero-Inflated Skellam distribution rows to the dist_table sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New rows 29-31, following the same pattern as the existing Skellam rows (20-22)
# and the other Zero-Inflated distribution rows (26-28).
# The G:H columns store the literal text "TRUE"/"FALSE" (not booleans), matching
# the existing rows above them, so a leading apostrophe forces text entry.

$ws.Range("A29").Value = "Zero-Inflated Skellam"
$ws.Range("B29").Value = "Difference"
$ws.Range("C29").Value = "ziskellam"
$ws.Range("D29").Value = "diff"
$ws.Range("E29").Value = "integer"
$ws.Range("F29").Value = "uni"
$ws.Range("G29").Value = "'FALSE"
$ws.Range("H29").Value = "'FALSE"

$ws.Range("A30").Value = "Zero-Inflated Skellam"
$ws.Range("B30").Value = "Mean-Dispersion"
$ws.Range("C30").Value = "ziskellam"
$ws.Range("D30").Value = "meandisp"
$ws.Range("E30").Value = "integer"
$ws.Range("F30").Value = "uni"
$ws.Range("G30").Value = "'FALSE"
$ws.Range("H30").Value = "'FALSE"

$ws.Range("A31").Value = "Zero-Inflated Skellam"
$ws.Range("B31").Value = "Mean-Variance"
$ws.Range("C31").Value = "ziskellam"
$ws.Range("D31").Value = "meanvar"
$ws.Range("E31").Value = "integer"
$ws.Range("F31").Value = "uni"
$ws.Range("G31").Value = "'FALSE"
$ws.Range("H31").Value = "'TRUE"

# Apply the same styles as the rows above (text format for G and H columns)
$ws.Range("G29:H31").Style = $ws.Range("G28:H28").Style

$ws.Range("B38").Select()
